$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(16124, 10961, 451, 189, 3512, 1011)
    3  = @(-0.1447299462922762, -0.1492223798200119, -0.07536776923588184, 0.2185821020819343, -0.2456433932824857, 0.06028432561899122)
    4  = @(-0.1485768202172113, -0.1549126428413103, -0.2220088286771384, -0.09631764782534602, -0.2651009412249949, 0.00720415154307652)
    5  = @(3.369940413048412, 3.386813866351882, 2.76339711883103, 2.701165729142938, 3.388435041038097, 3.479193489771451)
    6  = @(240, 165, 4, 5, 45, 21)
    7  = @(0.2968447583444246, 0.2850602023802188, 0.3847566437620272, 0.435241139923141, 0.2845470015691641, 0.3754272205442577)
    8  = @(0.29448179966315, 0.281520247394755, 0.3008598224568491, 0.2076517485488844, 0.2733712821837823, 0.3401480049683058)
    9  = @(2.641169716150813, 2.671308390905421, 2.090201999751256, 2.296365853398741, 2.5679871953144, 2.836428711636537)
    10 = @(339, 228, 7, 5, 65, 34)
    11 = @(0.4415747046367008, 0.4342825822002306, 0.4601244129979091, 0.2166590378412067, 0.5301903948516499, 0.3151428949252665)
    12 = @(0.4430586198803613, 0.4364328902360654, 0.5228686511339875, 0.3039693963742304, 0.5384722234087772, 0.3329438534252293)
    13 = @(0.7287706968975982, 0.7155054754464607, 0.6731951190797734, 0.4047998757441964, 0.820447845723697, 0.6427647781349139)
    14 = @(-99, -63, -3, 0, -20, -13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($col = 2; $col -le 7; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
